$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H13").Value = 4.1
$ws.Range("I13").Value = 5.75
$ws.Range("N13").Value = 1.7
$ws.Range("O13").Value = 2.1
$ws.Range("P13").Value = 1.3
$ws.Range("Q13").Value = 3.4
$ws.Range("U13").Value = 8
$ws.Range("AJ13").Value = 201

$ws.Range("G14").Value = 3.4
$ws.Range("H14").Value = 3.2
$ws.Range("K14").Value = 8.5
$ws.Range("Z14").Value = 8.5

$ws.Range("N19").Value = 1.9
$ws.Range("O19").Value = 1.9

$ws.Range("I20").Value = 3.5
$ws.Range("L20").Value = 1.4
$ws.Range("M20").Value = 2.75
$ws.Range("AE20").Value = 17

$ws.Range("J25").Value = 1.04
$ws.Range("K25").Value = 13

$ws.Range("G28").Value = 2.55
$ws.Range("I28").Value = 2.65
$ws.Range("Q28").Value = 2.25
$ws.Range("T28").Value = 6.6
$ws.Range("U28").Value = 11.25
$ws.Range("V28").Value = 10.5
$ws.Range("W28").Value = 28
$ws.Range("X28").Value = 26
$ws.Range("AC28").Value = 110
$ws.Range("AD28").Value = 6.7
$ws.Range("AE28").Value = 11.75
$ws.Range("AF28").Value = 10.75
$ws.Range("AG28").Value = 29
$ws.Range("AH28").Value = 27

$ws.Range("G31").Value = 4.2
$ws.Range("H31").Value = 3.5
$ws.Range("I31").Value = 1.91
$ws.Range("U31").Value = 21
$ws.Range("X31").Value = 34
$ws.Range("Y31").Value = 41
$ws.Range("Z31").Value = 10
$ws.Range("AD31").Value = 7.5
$ws.Range("AE31").Value = 9
$ws.Range("AF31").Value = 8.5

$ws.Range("R32").Value = 1.83
$ws.Range("AB32").Value = 15
$ws.Range("AC32").Value = 75

$wb.Save()
